$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.722.02"
$ws.Cells.Item(2, 5).Value = "  +0.53%  "
$ws.Cells.Item(3, 4).Value = "1.638.94"
$ws.Cells.Item(3, 5).Value = "  -0.50%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "212.77"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.14%  "
$ws.Cells.Item(6, 5).Value = "  -2.20%  "
$ws.Cells.Item(7, 5).Value = "  -0.02%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "23.24"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.58%  "
$ws.Cells.Item(9, 5).Value = "  +1.61%  "
$ws.Cells.Item(10, 5).Value = "  +0.01%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0888"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.02%  "
$ws.Cells.Item(12, 4).Value = "1.871.81"
$ws.Cells.Item(12, 5).Value = "  -0.44%  "
$ws.Cells.Item(13, 4).Value = "1.642.86"
$ws.Cells.Item(13, 5).Value = "  -0.13%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.05"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.44%  "
$ws.Cells.Item(15, 5).Value = "  -3.71%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "64.80"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.51%  "
$ws.Cells.Item(17, 4).Value = "27.682.40"
$ws.Cells.Item(17, 5).Value = "  +0.52%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "230.27"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.41%  "
$ws.Cells.Item(19, 5).Value = "  +1.90%  "
$ws.Cells.Item(20, 5).Value = "  -0.11%  "
$ws.Cells.Item(22, 5).Value = "  -0.24%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.23"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +4.91%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.09"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +3.78%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "151.03"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.71%  "
$ws.Cells.Item(27, 5).Value = "  -0.98%  "
$ws.Cells.Item(28, 5).Value = "  +0.06%  "
$ws.Cells.Item(29, 5).Value = "  +0.04%  "
$ws.Cells.Item(30, 5).Value = "  +0.27%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.0487"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.28%  "
$ws.Cells.Item(33, 4).Value = "1.463.92"
$ws.Cells.Item(33, 5).Value = "  +2.77%  "
$ws.Cells.Item(34, 5).Value = "  -2.08%  "
$ws.Cells.Item(35, 5).Value = "  -2.00%  "
$ws.Cells.Item(36, 5).Value = "  -0.42%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.568"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.06%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.883"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.31%  "
$ws.Cells.Item(39, 5).Value = "  -0.02%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.897"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +9.72%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "69.07"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +6.07%  "
$ws.Cells.Item(42, 5).Value = "  +0.02%  "
$ws.Cells.Item(43, 5).Value = "  -1.99%  "
$ws.Cells.Item(44, 5).Value = "  +1.36%  "
$ws.Cells.Item(45, 5).Value = "  -0.90%  "
$ws.Cells.Item(46, 5).Value = "  -0.66%  "
$ws.Cells.Item(47, 4).Value = "1.781.30"
$ws.Cells.Item(47, 5).Value = "  -0.48%  "
$ws.Cells.Item(48, 5).Value = "  +3.16%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "87.05"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.25%  "
$ws.Cells.Item(50, 5).Value = "  -1.22%  "
$ws.Cells.Item(51, 5).Value = "  -0.06%  "
